# Se permite comparar archivos XLSX
# Adds an "Adicional" column/cell (D5, underlined) to the sheet and
# normalizes the formatting of the existing data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the (already default) "Normal" style to the existing table so the
# workbook gets an explicit cell-format record for A1:C5, matching the
# re-saved file produced by the new XLSX tooling.
$ws.Range("A1:C5").Style = "Normal"

# New cell: D5 = "Adicional", underlined.
$ws.Range("D5").Value = "Adicional"
$ws.Range("D5").Font.Underline = $true

# Move the active selection to the newly added cell, like the source file.
$ws.Range("D5").Select() | Out-Null
